$d = $word.ActiveDocument

# --- Step 1: structural changes (bottom-up so earlier indices are stable) ---

# Merge paragraphs 24-26 ("כאשר a_i" / "" / "בנוסף השיפצורים") into a single paragraph
$d.Paragraphs.Item(26).Range.Delete()
$d.Paragraphs.Item(25).Range.Delete()

# Merge paragraphs 6-8 ("היום סוקרים" / "מנגנון תשומת לב" / "שכבת MLP") into a single paragraph
$d.Paragraphs.Item(8).Range.Delete()
$d.Paragraphs.Item(7).Range.Delete()

# Insert 3 new paragraphs after paragraph 4 (the paper-link URL paragraph)
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertParagraphAfter() | Out-Null
$p4.Range.InsertParagraphAfter() | Out-Null
$p4.Range.InsertParagraphAfter() | Out-Null

# --- Step 2: set final text for every paragraph, by position ---

$d.Paragraphs.Item(1).Range.Text = "Review 204: Simple linear attention language models balance the recall-throughput tradeoff"
$d.Paragraphs.Item(2).Range.Text = "Paper: https://arxiv.org/abs/2402.18668v2"
$d.Paragraphs.Item(4).Range.Text = "https://arxiv.org/abs/2402.18668"
$d.Paragraphs.Item(7).Range.Text = "מודלי שפה ענקיים של היום מפגינים יכולת מרשימה של למידת in-context כלומר יכולת לבצע משימות חדשות (שלא אומן עליהם באופן מפורש) בהתבסס על כמה דוגמאות המדגימות (ממחישות) את את המשימה. כמובן דוגמאות אלו מוזנות למודל שפה כפרומפט. המאמר שנסקור היום מדבר על משימת in-context ספציפית הנקראת recall. המטרה של משימה זו היא לזהות חוקיות מסומיות בפרומפט ולענות על שאלות בנוגע אליו. למשל אם פרומפט המוזן הוא `"A 4 B 3 C 6 F 1 G 2`". אם לאחר מכן אני מכניסים למודל שפה `"?B ? F ? C `" המודל צריך לענות 3 1 6 כלומר המספר בא מיד אחרי כל אות בפורמפט השאלה."
$d.Paragraphs.Item(9).Range.Text = "ארכיטקטורת הטרנספומרים מתמודדת בהצלחה עם משימות recall אך היא מתקשה עם אורכי הקשר (context length) מאוד ארוכים עקב מנגנון self-attention שלהם. ד״א המימושים המודרניים של מנגנון זה (כמו FlashAttention2 ו-Paged-Attention) הם בעלי סיבוכיות subquadratic במונחי אורך הסדרה אך עדיין גם הם מתקשים ״לעכל״ אורכי הקשר ממש ארוכים. "
$d.Paragraphs.Item(11).Range.Text = "כדי לתת מענה לסוגיה זו הוצעו מספר חלופות למנגנון ה-attention כמו attention לינארי, שיטות המבוססות על חלון הזז (sliding window) ובנוסף לאחרונה משפחת ארכיטקטורות ממבה (סקרתי אותן בהרחבה לפני כחודשיים). "
$d.Paragraphs.Item(13).Range.Text = "מנגנון attention לינארי בגדול מחליף את הסופטמקס של המכפלה הפנימית של וקטורי שאילתה (Q) ווקטורי ערך (K) למכפלה הפנימית של (f(Q ו- (f(K עבור פונקציה לא לינארית f (יש לא מעט מאמרים המציעים לקחת פונקציות f שונות עבור ההחלפה הזו). אחת הדוגמאות היא לבחור f בתור כמה איברים ראשונים של פיתוח טיילור של סופטמקס. "
$d.Paragraphs.Item(15).Range.Text = "פעולה זו מאפשרת להחליף סדר הפעולות בחישוב ה-attention ולבצע את החישוב באופן לינארי במונחי אורך הסדרה. דרך אגב החלפה זו היא כמו reparameterization trick ב- SVMs אבל בכיוון ההפוך. היא מאפשרת להיפטר מ״גרירה״ של הייצוגים של כל הטוקנים הקודמים באופן מפורש באינפרמס ומאפשרת חישוב בסגנון RNN. כלומר כל הזכרון עד טוקן n נדחס לכדי 2 וקטורים (ממליץ לקרוא על זה כאן) וכמובן זה מאפשר לחסוך במשאבי חישוב הנדרשים לביצוע אינפרנס באופן משמעותי. "
$d.Paragraphs.Item(17).Range.Text = "מנגנון ה-attention עם החלון הזז הוא פשוט הגבלת גודל ההקשר במנגנון ה-attention כאשר יש מגוון גישות ל`"איך לדחוס״ את הדאטה שלא נכנסת לחלון זה (העבר). בתוך החלון ה-attention מחושב באופן רגיל כלומר הגדלה משמעותית של חלון זה משפרת את הביצועים אבל גם כרוכה בביצוע של יותר חישובים."
$d.Paragraphs.Item(19).Range.Text = "מצד אחד ארכיטקטורות המבוססת על attention לינארי יודעות להסתדר לא רע עם אורכי הקשר ארוכים מאוד במשימות מסוימות אבל מתקשות לספק ביצועים גבוהים לשאלות בסגנון recall. מצד שני ארכיטקטורות הממשות חלון attention זז מסתדרות יפה עם משימות recall בתוך החלון הזה אולם כדי להביא ביצועים גבוהים עם הקשר ארוך צריך להגדיל את גודל החלון שכאמור כרוך בהקצאה של יותר משאבים ואו/גם ב-latencies גבוהים יותר באינפרנס. "
$d.Paragraphs.Item(21).Range.Text = "אוקיי דיברנו הרבה על הרקע למאמר אז הגיע הזמן לדבר על  מאמר עצמו. קודם כל החמברים מוכיחים באופן תיאורטי (את הקטע הזה הכי אהבתי כאן) כי ככל שאורך הקלט למשימת recall ״המודל צריך לזכור״ (O(N `"מידע`" כאשר N הוא ״אורך״ של פרומפט ה-recall (זה גם נבדק אמפירית). כלומר זה תקף לכל ארכיטקטורה והשאלה היחידה איך כל מודל (למשל טרנספורמר לינארי, hyena, mamba, s3 ועוד) בונים ומנהלים את הזכרון הזה ואיך הוא משפיע על ביצועי אינפרנס."
$d.Paragraphs.Item(23).Range.Text = "לגבי החידוש שהמאמר מציע: המחברים שילבו את ה״טוב״ שיש במנגנון ה-attention הלינארי ובגישת החלון הזז והציעו מנגנון attention חדש הנקרא Based. הם לקחו מנגנון ה-attention הלינארי החסכוני והיעיל מבחינת ניהול הזכרון והוסיפו לו חלון זז קצר יחסית המממש מנגנון attention רגיל של הטרנספורמים. וזה עבד להם לא רע בכלל במשימות recall שונות המצריכות חלון הקשר גדול. בנוסף גם הציעו מספר שכלולים לשיטה זו המאפשרים להריץ אותה בצורה מאוד יעילה על GPUs (למשל בחירת גודל החלון כדי שיהיה ניתן לבצע את החשובים עבור על ידי שימוש רק הזכרון המהיר של GPU."
$d.Paragraphs.Item(25).Range.Text = "בסך הכל מאמר די נחמד…"
